# Fixed #418 Empty AQL expressions generate empty lines.
#
# An AQL expression/field that resolves to an empty string used to leave a
# stray, completely empty paragraph behind in the generated document. Here
# that is the paragraph right after "Start of demonstration:" (it has no
# text at all - just an empty run) and right before "Some value". Remove it,
# which also removes its paragraph mark so the surrounding paragraphs become
# adjacent, exactly as if the empty line had never been generated.

$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    $text = $para.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark / cell mark;
    # an empty paragraph's text is just that mark, so strip it before
    # checking for emptiness.
    $trimmed = $text.Replace([char]13, "").Replace([char]7, "").Trim()
    if ($trimmed.Length -eq 0) {
        $para.Range.Delete()
    }
}
